$d = $word.ActiveDocument

# Word's Font.Color is a BGR-packed integer (wdColor), so turn the OOXML
# "345A8A" RGB hex value used by the template into that form.
$hex = "345A8A"
$r = [Convert]::ToInt32($hex.Substring(0,2), 16)
$g = [Convert]::ToInt32($hex.Substring(2,2), 16)
$b = [Convert]::ToInt32($hex.Substring(4,2), 16)
$abstractTitleColor = ($b * 65536) + ($g * 256) + $r

# 1. Add the new "Abstract Title" paragraph style (inserted just before "Abstract").
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true
$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0
$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = $abstractTitleColor

# 2. Tweak existing "Abstract" style: space-before goes from 300 (15pt) to 100 (5pt) twips.
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# 3. Add the new "Footnote Block Text" paragraph style (based on "Footnote Text").
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Host "Styles updated"
